$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 92: correct the date/time value in column A ---
$ws.Range("A92").Value = 45475.2916666667

# --- Add new row 93 ---

# Copy the formatting of row 92 onto row 93 first so the new row's cells
# (in particular A93's date number format / style) match the existing data
# rows exactly.
$ws.Range("A92:H92").Copy()
$ws.Range("A93:H93").PasteSpecial(-4122)

$ws.Range("A93").Value = 45476.5759375
$ws.Range("B93").Value = 9000
$ws.Range("C93").Value = 3.32999992370605
$ws.Range("D93").Value = 3.23000001907349
$ws.Range("E93").Value = 3.23000001907349
$ws.Range("F93").Value = 3.25999999046326

# Column G ("adj_close") stores the close price as text (shared string),
# matching the rest of the sheet, rather than as a number. Entering it as a
# formula that evaluates to the text and then collapsing the formula down to
# its value (paste values only) makes Excel store it as a genuine text cell
# without needing to change the cell's number format/style.
$ws.Range("G93").Formula = "=""3.25999999046326"""
$ws.Range("G93").Copy()
$ws.Range("G93").PasteSpecial(-4163)

$ws.Range("H93").Value = "ESPE.MI"
